$d = $word.ActiveDocument

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2, wdHeaderFooterEvenPages = 3
$wdHeaderFooterPrimary = 1
$wdHeaderFooterFirstPage = 2

# Rename each inline picture by round-tripping it through a floating Shape
# (the only object that exposes a writable .Name) and back to an InlineShape,
# which keeps the <wp:inline> markup intact while updating the docPr name.
function Rename-InlinePicture($range, $index, $newName) {
    $ishape = $range.InlineShapes.Item($index)
    $shape = $ishape.ConvertToShape()
    $shape.Name = $newName
    [void]$shape.ConvertToInlineShape()
}

foreach ($sec in $d.Sections) {

    # Header (first page) holds the BTec logo: image1.jpg -> image2.jpg
    $hdrFirst = $sec.Headers.Item($wdHeaderFooterFirstPage)
    if ($hdrFirst.Exists -and $hdrFirst.Range.InlineShapes.Count -gt 0) {
        Rename-InlinePicture $hdrFirst.Range 1 "image2.jpg"
    }

    # Footer (primary/default) holds a Pearson logo: image2.png -> image1.png
    $ftrPrimary = $sec.Footers.Item($wdHeaderFooterPrimary)
    if ($ftrPrimary.Exists -and $ftrPrimary.Range.InlineShapes.Count -gt 0) {
        Rename-InlinePicture $ftrPrimary.Range 1 "image1.png"
    }

    # Footer (first page) holds a Pearson logo too: image2.png -> image1.png
    $ftrFirst = $sec.Footers.Item($wdHeaderFooterFirstPage)
    if ($ftrFirst.Exists -and $ftrFirst.Range.InlineShapes.Count -gt 0) {
        Rename-InlinePicture $ftrFirst.Range 1 "image1.png"
    }
}
